# Crypto list refresh (GitHub Actions) - update Price/Volume(1h) columns,
# plus the Polkadot/WrappedEther row swap (rows 12-13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.273.60'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '1.679.63'
$ws.Range('E3').Value = '  +0.68%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.41'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5271'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2705'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.47%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06485'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.98'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07538'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.67%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.690.45'
$ws.Range('E12').Value = '  +1.05%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.525'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5803'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008513'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.62'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.60%  '
$ws.Range('D17').Value = '26.315.90'
$ws.Range('E17').Value = '  +0.37%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.923'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('E19').Value = '  +0.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.87'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '190.36'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.198'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.008'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '145.12'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.808'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1242'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.61%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.81'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06568'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.355'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.52%  '
$ws.Range('E30').Value = '  +0.82%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.597'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.587'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.97%  '
$ws.Range('E33').Value = '  +1.32%  '
$ws.Range('E34').Value = '  +2.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6231'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.63%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.403'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.31%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.737'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.444'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.73%  '
$ws.Range('D39').Value = '1.112.40'
$ws.Range('E39').Value = '  +2.75%  '
$ws.Range('E40').Value = '  +1.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8755'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.014'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.42%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.82'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.29%  '
$ws.Range('D44').Value = '1.829.99'
$ws.Range('E44').Value = '  +0.72%  '
$ws.Range('E45').Value = '  +0.66%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.97'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.152'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.88%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.006'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05272'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.103'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4293'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.03%  '
